$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 121.375
$ws.Range("I33").Value = 114.23077
$ws.Range("K33").Value = 114.23077
$ws.Range("M33").Value = 114.76923

$ws.Range("H132").Value = 1606.4193
$ws.Range("I132").Value = 1626.6333
$ws.Range("K132").Value = 4879.8999
$ws.Range("M132").Value = -2349.8999

$ws.Range("H133").Value = 76107.164
$ws.Range("J133").Value = 76107.164
$ws.Range("L133").Value = 76107.164
$ws.Range("N133").Value = -86227.164

$ws.Range("H134").Value = 39996.273
$ws.Range("J134").Value = 39996.273
$ws.Range("L134").Value = 39996.273
$ws.Range("N134").Value = -50136.273

$ws.Range("H135").Value = 2785.4
$ws.Range("I135").Value = 2856.75
$ws.Range("J135").Value = 2500
$ws.Range("K135").Value = 25710.75
$ws.Range("L135").Value = 22500
$ws.Range("M135").Value = -23175.75
$ws.Range("N135").Value = -27570

$ws.Range("H136").Value = 57579.5
$ws.Range("J136").Value = 57579.5
$ws.Range("L136").Value = 57579.5
$ws.Range("N136").Value = -67779.5

$ws.Range("H137").Value = 606905.4
$ws.Range("I137").Value = 1795.9375
$ws.Range("J137").Value = 1817124.2
$ws.Range("K137").Value = 5387.8125
$ws.Range("L137").Value = 5451372.6
$ws.Range("M137").Value = -2837.8125
$ws.Range("N137").Value = -5456472.6

$ws.Range("H139").Value = 74227.5
$ws.Range("J139").Value = 74227.5
$ws.Range("L139").Value = 74227.5
$ws.Range("N139").Value = -84507.5

$ws.Range("H140").Value = 63871.918
$ws.Range("J140").Value = 65068.547
$ws.Range("L140").Value = 65068.547
$ws.Range("N140").Value = -75428.54699999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9877.359
$ws.Range("I32").Value = 3877.4783
$ws.Range("K32").Value = 3877.4783
$ws.Range("M32").Value = -3590.4783

$ws.Range("H61").Value = 33254.094
$ws.Range("I61").Value = 1825.5186
$ws.Range("K61").Value = 1825.5186
$ws.Range("M61").Value = -1613.5186

$ws.Range("H136").Value = 33254.094
$ws.Range("I136").Value = 1825.5186
$ws.Range("K136").Value = 5476.5558
$ws.Range("M136").Value = -2926.5558

$ws.Range("H139").Value = 124857.5
$ws.Range("J139").Value = 124857.5
$ws.Range("L139").Value = 124857.5
$ws.Range("N139").Value = -135137.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 777
$ws.Range("I5").Value = 647.8
$ws.Range("J5").Value = 1100
$ws.Range("K5").Value = 647.8
$ws.Range("L5").Value = 1100
$ws.Range("M5").Value = -534.8
$ws.Range("N5").Value = -1326

$ws.Range("H132").Value = 33047.145
$ws.Range("J132").Value = 33047.145
$ws.Range("L132").Value = 33047.145
$ws.Range("N132").Value = -43167.145

$ws.Range("H135").Value = 48870.57
$ws.Range("J135").Value = 48870.57
$ws.Range("L135").Value = 48870.57
$ws.Range("N135").Value = -59010.57

$ws.Range("H138").Value = 72985.586
$ws.Range("J138").Value = 72985.586
$ws.Range("L138").Value = 72985.586
$ws.Range("N138").Value = -83265.586

$ws.Range("H140").Value = 74796.8
$ws.Range("J140").Value = 74796.8
$ws.Range("L140").Value = 74796.8
$ws.Range("N140").Value = -85156.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1651.0588
$ws.Range("J16").Value = 2376.125
$ws.Range("L16").Value = 2376.125
$ws.Range("N16").Value = -2950.125

$ws.Range("H58").Value = 1733.6666
$ws.Range("I58").Value = 1595.8334
$ws.Range("K58").Value = 1595.8334
$ws.Range("M58").Value = -1392.8334

$ws.Range("H93").Value = 63250
$ws.Range("I93").Value = 63250
$ws.Range("K93").Value = 63250
$ws.Range("M93").Value = -61378

$ws.Range("H113").Value = 1651.0588
$ws.Range("J113").Value = 2376.125
$ws.Range("L113").Value = 2376.125
$ws.Range("N113").Value = -6716.125

$ws.Range("H136").Value = 1733.6666
$ws.Range("I136").Value = 1595.8334
$ws.Range("K136").Value = 4787.5002
$ws.Range("M136").Value = -2237.5002

$ws.Range("H138").Value = 99996
$ws.Range("J138").Value = 99996
$ws.Range("L138").Value = 99996
$ws.Range("N138").Value = -110276

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 17153886
$ws.Range("I4").Value = 10000046
$ws.Range("J4").Value = 41000020
$ws.Range("K4").Value = 30000138
$ws.Range("L4").Value = 123000060
$ws.Range("M4").Value = -30000026
$ws.Range("N4").Value = -123000284

$ws.Range("H16").Value = 10.5
$ws.Range("I16").Value = 10.5
$ws.Range("K16").Value = 31.5
$ws.Range("M16").Value = 141.5

$ws.Range("H115").Value = 6464
$ws.Range("I115").Value = 928
$ws.Range("J115").Value = 12000
$ws.Range("K115").Value = 2784
$ws.Range("L115").Value = 36000
$ws.Range("M115").Value = -1609
$ws.Range("N115").Value = -38350

$ws.Range("H122").Value = 505418.56
$ws.Range("J122").Value = 721833.6
$ws.Range("L122").Value = 6496502.399999999
$ws.Range("N122").Value = -6501402.399999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 311.42856
$ws.Range("I2").Value = 276.125
$ws.Range("K2").Value = 276.125
$ws.Range("M2").Value = -163.125

$ws.Range("H19").Value = 15174
$ws.Range("I19").Value = 13243.6
$ws.Range("J19").Value = 20000
$ws.Range("K19").Value = 13243.6
$ws.Range("L19").Value = 20000
$ws.Range("M19").Value = -12955.6
$ws.Range("N19").Value = -20576

$ws.Range("H93").Value = 12717.556
$ws.Range("J93").Value = 12717.556
$ws.Range("L93").Value = 12717.556
$ws.Range("N93").Value = -16461.556

$ws.Range("H97").Value = 347.41666
$ws.Range("I97").Value = 261.75
$ws.Range("J97").Value = 518.75
$ws.Range("K97").Value = 261.75
$ws.Range("L97").Value = 518.75
$ws.Range("M97").Value = 234.25
$ws.Range("N97").Value = -1510.75

$ws.Range("H135").Value = 45436.25
$ws.Range("J135").Value = 45436.25
$ws.Range("L135").Value = 45436.25
$ws.Range("N135").Value = -55576.25

$ws.Range("H140").Value = 95181.664
$ws.Range("J140").Value = 95607.27
$ws.Range("L140").Value = 95607.27
$ws.Range("N140").Value = -105967.27

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4915.6924
$ws.Range("J7").Value = 5174
$ws.Range("L7").Value = 5174
$ws.Range("N7").Value = -5398

$ws.Range("H55").Value = 4189.7
$ws.Range("I55").Value = 1309
$ws.Range("J55").Value = 6796.048
$ws.Range("K55").Value = 1309
$ws.Range("L55").Value = 6796.048
$ws.Range("M55").Value = -1136
$ws.Range("N55").Value = -7142.048

$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

$ws.Range("H122").Value = 30004400
$ws.Range("I122").Value = 40003292
$ws.Range("J122").Value = 13339580
$ws.Range("K122").Value = 120009876
$ws.Range("L122").Value = 40018740
$ws.Range("M122").Value = -120007426
$ws.Range("N122").Value = -40023640

$ws.Range("H126").Value = 4915.6924
$ws.Range("J126").Value = 5174
$ws.Range("L126").Value = 15522
$ws.Range("N126").Value = -20462

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()

$ws.Range("H100").Value = 3247558.5
$ws.Range("I100").Value = 3969064
$ws.Range("J100").Value = 783.25
$ws.Range("K100").Value = 7938128
$ws.Range("L100").Value = 1566.5
$ws.Range("M100").Value = -7937587
$ws.Range("N100").Value = -2648.5

$ws.Range("H138").Value = 150214.5
$ws.Range("J138").Value = 150214.5
$ws.Range("L138").Value = 150214.5
$ws.Range("N138").Value = -160494.5
